# This script applies the numeric updates described in the commit diff
# to the "ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR" market-data sheets.
# All target cells are plain numeric values (no formulas in the workbook),
# so each change is applied as a direct Range(...).Value assignment.
# A couple of cells are removed entirely (cleared) or newly introduced,
# matching the cell-level add/remove/change operations in the diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 42010.5
$ws.Range("I34").Value = 42010.5
$ws.Range("K34").Value = 42010.5
$ws.Range("M34").Value = -41739.5
$ws.Range("H40").Value = 32299.666
$ws.Range("I40").Value = 29999
$ws.Range("J40").Value = 33450
$ws.Range("K40").Value = 29999
$ws.Range("L40").Value = 33450
$ws.Range("M40").Value = -29823
$ws.Range("N40").Value = -33802
$ws.Range("H45").Value = 1337.091
$ws.Range("I45").Value = 1329.3334
$ws.Range("K45").Value = 1329.3334
$ws.Range("M45").Value = -952.3334
$ws.Range("H74").Value = 2316827
$ws.Range("I74").Value = 1030991.7
$ws.Range("K74").Value = 1030991.7
$ws.Range("M74").Value = -1030117.7
$ws.Range("H77").Value = 2316827
$ws.Range("I77").Value = 1030991.7
$ws.Range("K77").Value = 5154958.5
$ws.Range("M77").Value = -5150590.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 401.89474
$ws.Range("I22").Value = 299.55554
$ws.Range("J22").Value = 494
$ws.Range("K22").Value = 299.55554
$ws.Range("L22").Value = 494
$ws.Range("M22").Value = -126.55554
$ws.Range("N22").Value = -840
$ws.Range("H86").Value = 2398.5557
$ws.Range("I86").Value = 2062.2144
$ws.Range("J86").Value = 3575.75
$ws.Range("K86").Value = 2062.2144
$ws.Range("L86").Value = 3575.75
$ws.Range("M86").Value = -939.2143999999998
$ws.Range("N86").Value = -5821.75
$ws.Range("H89").Value = 2398.5557
$ws.Range("I89").Value = 2062.2144
$ws.Range("J89").Value = 3575.75
$ws.Range("K89").Value = 10311.072
$ws.Range("L89").Value = 17878.75
$ws.Range("M89").Value = -4695.072
$ws.Range("N89").Value = -29110.75
$ws.Range("H105").Value = 6033.4
$ws.Range("I105").Value = 6037.1113
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 6037.1113
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -4290.1113
$ws.Range("N105").Value = -9494
$ws.Range("H134").Value = 15354293
$ws.Range("I134").Value = 7146508
$ws.Range("K134").Value = 21439524
$ws.Range("M134").Value = -21436989

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19152.861
$ws.Range("I31").Value = 34563
$ws.Range("J31").Value = 2642
$ws.Range("K31").Value = 34563
$ws.Range("L31").Value = 2642
$ws.Range("M31").Value = -34268
$ws.Range("N31").Value = -3232
$ws.Range("H34").Value = 19152.861
$ws.Range("I34").Value = 34563
$ws.Range("J34").Value = 2642
$ws.Range("K34").Value = 34563
$ws.Range("L34").Value = 2642
$ws.Range("M34").Value = -34361
$ws.Range("N34").Value = -3046
$ws.Range("H62").Value = 35974.25
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 35974.25
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 35974.25
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -37222.25
$ws.Range("H65").Value = 35974.25
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 35974.25
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 179871.25
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -186111.25
$ws.Range("H99").Value = 2175.3333
$ws.Range("I99").Value = 2006
$ws.Range("K99").Value = 2006
$ws.Range("M99").Value = -508
$ws.Range("H122").Value = 2025.6666
$ws.Range("I122").Value = 1433
$ws.Range("K122").Value = 4299
$ws.Range("M122").Value = -1849
$ws.Range("H126").Value = 2175.3333
$ws.Range("I126").Value = 2006
$ws.Range("K126").Value = 6018
$ws.Range("M126").Value = -3548
$ws.Range("H134").Value = 10002545
$ws.Range("I134").Value = 2820
$ws.Range("K134").Value = 8460
$ws.Range("M134").Value = -5925

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 844.8889
$ws.Range("I8").Value = 844.8889
$ws.Range("K8").Value = 2534.6667
$ws.Range("M8").Value = -2395.6667
$ws.Range("H40").Value = 141.5
$ws.Range("J40").Value = 153.6923
$ws.Range("L40").Value = 614.7692
$ws.Range("N40").Value = -752.7692
$ws.Range("H122").Value = 85854.836
$ws.Range("I122").Value = 574
$ws.Range("J122").Value = 102911
$ws.Range("K122").Value = 5166
$ws.Range("L122").Value = 926199
$ws.Range("M122").Value = -2716
$ws.Range("N122").Value = -931099

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4584.615
$ws.Range("J80").Value = 5685
$ws.Range("L80").Value = 5685
$ws.Range("N80").Value = -7681
$ws.Range("H83").Value = 4584.615
$ws.Range("J83").Value = 5685
$ws.Range("L83").Value = 28425
$ws.Range("N83").Value = -38409
$ws.Range("H113").Value = 2414.6667
$ws.Range("J113").Value = 2506
$ws.Range("L113").Value = 2506
$ws.Range("N113").Value = -6846
$ws.Range("H126").Value = 4842.2856
$ws.Range("I126").Value = 7115.5
$ws.Range("K126").Value = 21346.5
$ws.Range("M126").Value = -18876.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2383.3
$ws.Range("I68").Value = 2380.3333
$ws.Range("K68").Value = 2380.3333
$ws.Range("M68").Value = -1631.3333
$ws.Range("H71").Value = 2383.3
$ws.Range("I71").Value = 2380.3333
$ws.Range("K71").Value = 11901.6665
$ws.Range("M71").Value = -8157.666499999999
$ws.Range("H82").Value = 2595.3333
$ws.Range("I82").Value = 979.6
$ws.Range("J82").Value = 4615
$ws.Range("K82").Value = 979.6
$ws.Range("L82").Value = 4615
$ws.Range("M82").Value = -618.6
$ws.Range("N82").Value = -5337
$ws.Range("H85").Value = 2595.3333
$ws.Range("I85").Value = 979.6
$ws.Range("J85").Value = 4615
$ws.Range("K85").Value = 979.6
$ws.Range("L85").Value = 4615
$ws.Range("M85").Value = 268.4
$ws.Range("N85").Value = -7111
$ws.Range("H122").Value = 3222.3635
$ws.Range("I122").Value = 2850
$ws.Range("K122").Value = 8550
$ws.Range("M122").Value = -6100
$ws.Range("H136").Value = 3161
$ws.Range("I136").Value = 2526.8572
$ws.Range("J136").Value = 7600
$ws.Range("K136").Value = 7580.571599999999
$ws.Range("L136").Value = 22800
$ws.Range("M136").Value = -5030.571599999999
$ws.Range("N136").Value = -27900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 18382.334
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 18382.334
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 18382.334
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -19138.334
$ws.Range("H62").Value = 3350.125
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3350.125
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740
$ws.Range("H96").Value = 3466.8
$ws.Range("I96").Value = 7500.5
$ws.Range("K96").Value = 7500.5
$ws.Range("M96").Value = -6127.5
$ws.Range("H126").Value = 2949.25
$ws.Range("I126").Value = 3099
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 9297
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -6827
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 2331.8064
$ws.Range("I132").Value = 2359.5334
$ws.Range("K132").Value = 7078.600199999999
$ws.Range("M132").Value = -4548.600199999999
$ws.Range("H136").Value = 2115.5833
$ws.Range("I136").Value = 2035.2727
$ws.Range("K136").Value = 6105.8181
$ws.Range("M136").Value = -3555.8181

